$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function FindParaRange([string]$text) {
    $r = $d.Content
    $ok = $r.Find.Execute($text, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Anchor text not found: $text" }
    return $r.Paragraphs(1).Range
}

# --- 1. Remove the stray _GoBack bookmark that currently sits after
#        "...approximately 1.5 times longer than it took with Eggplant." ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Collapse the whole "Scenarion Creation Log" / time-log block down to
#        just the (now relocated) "Time to Execute Common Scenario" heading. ---
$blockStart = FindParaRange("Scenarion Creation Log:")
$blockEnd   = FindParaRange("TODO: Check these numbers!!!")
$fullBlock  = $d.Range($blockStart.Start, $blockEnd.End)
$fullBlock.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="100" w:name="_Toc380582923"/><w:bookmarkStart w:id="101" w:name="_Toc381349842"/><w:bookmarkStart w:id="102" w:name="_Toc381349932"/><w:r><w:t>Time to Execute Common Scenario</w:t></w:r><w:bookmarkEnd w:id="100"/><w:bookmarkEnd w:id="101"/><w:bookmarkEnd w:id="102"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 3. "2 minutes 45 seconds..." -> "2 minutes 20 seconds..." ---
$bullet1 = FindParaRange("2 minutes 45 seconds for successful completion of the entire scenario.")
$bullet1.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>2 minutes 20</w:t></w:r><w:r><w:t xml:space="preserve"> seconds for successful completion of the entire scenario.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 4. Merge the "1 minute 20 seconds" bullet into a single run ---
$bullet2 = FindParaRange("1 minute 20 seconds when failure was induced.")
$bullet2.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>1 minute 20 seconds when failure was induced.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 5. Extend the Eggplant comparison sentence with the ATRT comparison ---
$eggplant = FindParaRange("These times are about 50% slower than Eggplant execution times.")
$eggplant.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>These times are about 50% slower than Eggplant execution times</w:t></w:r><w:r><w:t xml:space="preserve">, and slightly faster than ATRT execution times.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 6. Drop the "ADD BETTER EVAL HERE WHEN ATRT COMPLETE??" paragraph entirely ---
$addEval = FindParaRange("ADD BETTER EVAL HERE WHEN ATRT COMPLETE??")
$addEval.Delete()

# --- 7. "SUT Interaction and Performance" heading: renumber bookmarks, drop the
#        page-break run, and park the new _GoBack bookmark on this heading ---
$sut = FindParaRange("SUT Interaction and Performance")
$sut.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:bookmarkStart w:id="100" w:name="_Toc380582924"/><w:bookmarkStart w:id="101" w:name="_Toc381349843"/><w:bookmarkStart w:id="102" w:name="_Toc381349933"/><w:bookmarkStart w:id="103" w:name="_GoBack"/><w:bookmarkEnd w:id="103"/><w:r><w:t xml:space="preserve">SUT Interaction </w:t></w:r><w:r><w:t>and Performance</w:t></w:r><w:bookmarkEnd w:id="100"/><w:bookmarkEnd w:id="101"/><w:bookmarkEnd w:id="102"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 8. "Optical Character Recognition" heading now starts a fresh page ---
$ocr = FindParaRange("Optical Character Recognition")
$ocr.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:bookmarkStart w:id="100" w:name="_Toc380582926"/><w:bookmarkStart w:id="101" w:name="_Toc381349845"/><w:bookmarkStart w:id="102" w:name="_Toc381349935"/><w:r><w:lastRenderedPageBreak/><w:t>Optical Character Recognition</w:t></w:r><w:bookmarkEnd w:id="100"/><w:bookmarkEnd w:id="101"/><w:bookmarkEnd w:id="102"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 9. Merge the "debug code" runs and drop the mid-sentence page break ---
$dbg = FindParaRange("likely requires")
$dbg.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Sikuli’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> default IDE is the only one evaluated that does not have a debugger.  Debugging in the default IDE is possible but it relies on the creativi</w:t></w:r><w:r><w:t xml:space="preserve">ty of the developer.  In most cases it </w:t></w:r><w:r><w:t xml:space="preserve">likely requires “debug” code to be placed at strategic places within the source to allow the executing script to stop as desired.  At best it is cumbersome.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)

# --- 10. Start the "Use Sikuli to manipulate..." bullet with a page break ---
$useBullet = FindParaRange("manipulate the")
$useBullet.InsertXML(@'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>U</w:t></w:r><w:r><w:t xml:space="preserve">se </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sikuli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to manipulate the </w:t></w:r><w:r><w:t>SUT</w:t></w:r><w:r><w:t xml:space="preserve"> (r</w:t></w:r><w:r><w:t xml:space="preserve">hel10) via the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vncviewer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> window.  </w:t></w:r><w:r><w:t xml:space="preserve">NOTE: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:r><w:t>ikuli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> IDE displays outside of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vnc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - we found that this makes </w:t></w:r><w:r><w:t>it easier to capture screen im</w:t></w:r><w:r><w:t xml:space="preserve">ages for menus this way because </w:t></w:r><w:r><w:t xml:space="preserve">the menus stay open in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vncvie</w:t></w:r><w:r><w:t>wer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> when we activate the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:r><w:t>ikuli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>IDE for capture.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@)
